# --- Update project report and presentation ---

$wb = $excel.ActiveWorkbook

# 1) Hide the "GPT-4.1-mini" sheet (state="hidden").
$wsGpt   = $wb.Worksheets.Item("GPT-4.1-mini")
$wsMed   = $wb.Worksheets.Item("MedGemma 4B")
$wsDpo   = $wb.Worksheets.Item("MedGemma 4B DPO")

# 2) Fix up the "GPT-4.1-mini" sheet's selection before hiding it.
$wsGpt.Range("B8").Select()
$wsGpt.Visible = $false

# 3) "MedGemma 4B" sheet: zoom 70%, refreshed pane/selection, and convert the
#    row-12 AVERAGE formulas into one filled/shared formula (D12:K12).
$wsMed.Activate()
$excel.ActiveWindow.Zoom = 70
$wsMed.Range("D12:K12").Formula = "=AVERAGE(D2:D11)"
$wsMed.Range("G3").Select()

# 4) "MedGemma 4B DPO" sheet: rename the chat-log references from the
#    "it" run to the new "dpo" run, bump the Search-action score for case 6,
#    zoom 85%, and refresh the pane/selection. This sheet stays the active
#    tab (it was already the active tab in the source file).
$wsDpo.Activate()
$wsDpo.Range("L2").Value  = "logs\medgemma-4b-dpo-2025_07_20-02_34.txt"
$wsDpo.Range("L3").Value  = "logs\medgemma-4b-dpo-2025_07_20-02_46.txt"
$wsDpo.Range("L4").Value  = "logs\medgemma-4b-dpo-2025_07_20-22_52.txt"
$wsDpo.Range("L5").Value  = "logs\medgemma-4b-dpo-2025_07_20-23_05.txt"
$wsDpo.Range("L6").Value  = "logs\medgemma-4b-dpo-2025_07_21-01_36.txt"
$wsDpo.Range("L7").Value  = "logs\medgemma-4b-dpo-2025_07_21-01_42.txt"
$wsDpo.Range("L9").Value  = "logs\medgemma-4b-dpo-2025_07_21-02_59.txt"
$wsDpo.Range("L10").Value = "logs\medgemma-4b-dpo-2025_07_21-03_29.txt"
$wsDpo.Range("L11").Value = "logs\medgemma-4b-dpo-2025_07_21-03_50.txt"

$wsDpo.Range("G7").Value = 1

$excel.ActiveWindow.Zoom = 85
$wsDpo.Range("L2").Select()

Write-Output "edit applied"
